$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "28.115.39"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  +1.03%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.817.03"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  +2.59%  "

$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  -0.64%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "338.42"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  -0.29%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.9975"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  -0.81%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4382"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  +14.63%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3523"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  +3.50%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "45.65"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  -1.10%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.154"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  +1.47%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07425"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  +0.55%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "22.93"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  -1.47%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.9984"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  -0.84%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.264"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  -1.44%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "1.819.10"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  +2.48%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "7.290"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  -2.51%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.00001085"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  +1.21%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.06696"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  +0.21%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "82.19"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  +0.10%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.9992"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  -0.57%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.481"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  +0.72%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "17.28"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  -0.42%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "28.176.50"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  +1.18%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "12.03"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  +0.07%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.373"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  -1.07%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.497"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  +3.71%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "20.78"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  +0.71%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "155.67"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  +1.33%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.025.44"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  +2.46%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.304"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  -11.72%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "132.82"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  -0.72%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.048"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  +0.09%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.970"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  -0.64%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.09345"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  +5.21%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "12.35"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  -2.06%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.02371"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  -0.84%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.6774"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  -0.48%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "5.249"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  -0.98%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.06242"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  -2.02%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.2163"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  +0.32%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.481"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  -1.34%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.219"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  -0.71%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "8.221"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  +0.40%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.9979"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  -0.69%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "14.05"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  +0.01%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "3.878"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  +0.27%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.6141"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  -1.13%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "129.32"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  -3.21%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.046"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  -0.78%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.177"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  -2.18%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.07110"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  -3.60%  "
